$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date column cells remain plain text (not auto-converted to Excel date serials)
$dateRange = $ws.Range("A2:A18")
$dateRange.NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2024-09-05"
$rowVals = @(1, 19189, 5219, 8740, 0, 5529, 0, 15018, 19463, 73158, 19190.9189, 5219.5219, 8740.874, 0, 5529.5529, 0, 15019.5018, 19464.9463, 73165.3158)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(2, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(3, 1).Value = "2024-09-05"
$rowVals = @(2, 18379, 4985, 8351, 0, 5430, 0, 14349, 18477, 69971, 18380.8379, 4985.4985, 8351.8351, 0, 5430.543, 0, 14350.4349, 18478.8477, 69977.9971)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(3, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(4, 1).Value = "2024-09-05"
$rowVals = @(3, 17470, 4704, 8155, 0, 5207, 0, 13877, 17820, 67233, 17471.747, 4704.4704, 8155.8155, 0, 5207.5207, 0, 13878.3877, 17821.782, 67239.7233)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(4, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(5, 1).Value = "2024-09-05"
$rowVals = @(4, 17438, 4638, 8102, 0, 5178, 0, 13798, 17620, 66774, 17439.7438, 4638.4638, 8102.8102, 0, 5178.5178, 0, 13799.3798, 17621.762, 66780.6774)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(5, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(6, 1).Value = "2024-09-05"
$rowVals = @(5, 18041, 4882, 8235, 0, 5226, 0, 14143, 17900, 68427, 18042.8041, 4882.4882, 8235.8235, 0, 5226.5226, 0, 14144.4143, 17901.79, 68433.8427)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(6, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(7, 1).Value = "2024-09-05"
$rowVals = @(6, 18811, 5371, 8886, 0, 5719, 0, 14910, 19136, 72833, 18812.8811, 5371.5371, 8886.8886, 0, 5719.5719, 0, 14911.491, 19137.9136, 72840.2833)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(7, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(8, 1).Value = "2024-09-05"
$rowVals = @(7, 18689, 4586, 9006, 0, 7021, 0, 14958, 18598, 72858, 18690.8689, 4586.4586, 9006.9006, 0, 7021.7021, 0, 14959.4958, 18599.8598, 72865.2858)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(8, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(9, 1).Value = "2024-09-05"
$rowVals = @(8, 21111, 4741, 10697, 0, 10465, 0, 17456, 20215, 84685, 21113.1111, 4741.4741, 10698.0697, 0, 10466.0465, 0, 17457.7456, 20217.0215, 84693.4685)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(9, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(10, 1).Value = "2024-09-05"
$rowVals = @(9, 27401, 5613, 13587, 0, 12996, 0, 19031, 23924, 102552, 27403.7401, 5613.5613, 13588.3587, 0, 12997.2996, 0, 19032.9031, 23926.3924, 102562.2552)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(10, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(11, 1).Value = "2024-09-05"
$rowVals = @(10, 29301, 5688, 17571, 0, 14723, 0, 20103, 24871, 112257, 29303.9301, 5688.5688, 17572.7571, 0, 14724.4723, 0, 20105.0103, 24873.4871, 112268.2257)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(11, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(12, 1).Value = "2024-09-05"
$rowVals = @(11, 30596, 5902, 18831, 0, 14864, 0, 20562, 25402, 116157, 30599.0596, 5902.5902, 18832.8831, 0, 14865.4864, 0, 20564.0562, 25404.5402, 116168.6157)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(12, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(13, 1).Value = "2024-09-05"
$rowVals = @(12, 30033, 6037, 18926, 0, 14319, 0, 20138, 25728, 115181, 30036.0033, 6037.6037, 18927.8926, 0, 14320.4319, 0, 20140.0138, 25730.5728, 115192.5181)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(13, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(14, 1).Value = "2024-09-05"
$rowVals = @(13, 30101, 5877, 19038, 0, 15795, 0, 21958, 25669, 118438, 30104.0101, 5877.5877, 19039.9038, 0, 15796.5795, 0, 21960.1958, 25671.5669, 118449.8438)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(14, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(15, 1).Value = "2024-09-05"
$rowVals = @(14, 31220, 6127, 19084, 0, 15071, 0, 21170, 26296, 118968, 31223.122, 6127.6127, 19085.9084, 0, 15072.5071, 0, 21172.117, 26298.6296, 118979.8968)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(15, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(16, 1).Value = "2024-09-05"
$rowVals = @(15, 32197, 6332, 19180, 0, 15521, 0, 21620, 26956, 121806, 32200.2197, 6332.6332, 19181.918, 0, 15522.5521, 0, 21622.162, 26958.6956, 121818.1806)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(16, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(17, 1).Value = "2024-09-05"
$rowVals = @(16, 29317, 5863, 18637, 0, 14787, 0, 20501, 25861, 114966, 29319.9317, 5863.5863, 18638.8637, 0, 14788.4787, 0, 20503.0501, 25863.5861, 114977.4966)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(17, $c).Value = $rowVals[$c - 2] }

$ws.Cells.Item(18, 1).Value = "2024-09-05"
$rowVals = @(17, 0, 5726, 0, 0, 0, 0, 0, 9830, 15556, 0, 5726.5726, 0, 0, 0, 0, 0, 9830.983, 15557.5556)
for ($c = 2; $c -le 20; $c++) { $ws.Cells.Item(18, $c).Value = $rowVals[$c - 2] }

# Restore default style on date column so no unintended style/format change is introduced
$dateRange.Style = "Normal"
